# Outstandings.xlsx — "Add files via upload" edit
#
# Summary of the change being applied:
#  - Sheet "Purchase 22-23" (sheet1): a new ledger line is inserted as row 29
#    (payment dated 45197, ref "54/23-24", vendor "Namrata Rubber Product Pvt
#    Ltd", amount 26491), pushing the old rows 30-33 down to 31-34. The
#    running-total formulas renumber automatically; the "total so far" bold
#    boxed style moves from the old F28 onto the new F29, while F28 reverts to
#    the plain boxed style used by the rest of the block.
#  - Sheet "Sale 22-23" (sheet2): row 8's reference changes from the old
#    "b22-23MQ208" text to "b23-24MQ208", and a new row 9 is inserted (date
#    45199, ref "b23-24MQ210", vendor "Putzmeister Concrete Machines Pvt
#    Ltd", amount 39747) which now carries the running-total formula that used
#    to live on F8 (F8 itself becomes a plain, formula-less cell).
#  - The active tab switches from "Sale 22-23" back to "Purchase 22-23", and
#    each sheet's remembered selection moves to a new cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Purchase 22-23"
$ws2 = $wb.Worksheets.Item(2)   # "Sale 22-23"

# ---------------------------------------------------------------------------
# 1) Sheet1 ("Purchase 22-23"): insert the new row 29
# ---------------------------------------------------------------------------

$ws1.Rows("29").Insert()

# Bring over the plain "data row" look (boxed border, regular weight) from
# row 27 for columns A:E of the freshly inserted row.
$ws1.Range("A27:E27").Copy()
$ws1.Range("A29:E29").PasteSpecial(-4122)   # xlPasteFormats

# The bold "running total" box that used to sit on F28 moves down onto F29.
$ws1.Range("F28").Copy()
$ws1.Range("F29").PasteSpecial(-4122)       # xlPasteFormats

# ... and F28 itself drops back to the plain (non-bold) boxed style, matching
# the other rows in the block (same look as F27).
$ws1.Range("F27").Copy()
$ws1.Range("F28").PasteSpecial(-4122)       # xlPasteFormats

$ws1.Application.CutCopyMode = $false

$ws1.Rows("29").RowHeight = 14.4

$ws1.Range("A29").Value = $null
$ws1.Range("B29").Value = 45197
$ws1.Range("C29").Value = "54/23-24"
$ws1.Range("D29").Value = "Namrata Rubber Product Pvt Ltd"
$ws1.Range("E29").Value = 26491
$ws1.Range("F29").Formula = "=F28+E29"

# ---------------------------------------------------------------------------
# 2) Sheet2 ("Sale 22-23"): row 8's reference text + new row 9
# ---------------------------------------------------------------------------

$ws2.Rows("9").Insert()

# Row 9 takes on the same boxed look as row 8 (and the rest of that block).
$ws2.Range("A8:F8").Copy()
$ws2.Range("A9:F9").PasteSpecial(-4122)     # xlPasteFormats
$ws2.Application.CutCopyMode = $false

$ws2.Range("C8").Value = "b23-24MQ208"

$ws2.Range("A9").Value = $null
$ws2.Range("B9").Value = 45199
$ws2.Range("C9").Value = "b23-24MQ210"
$ws2.Range("D9").Value = "Putzmeister Concrete Machines Pvt Ltd"
$ws2.Range("E9").Value = 39747
$ws2.Range("F9").Formula = "=E5+E6+E7+E8+E9"

# F8 no longer carries the running-total formula (that moved to F9); it
# becomes a blank, formula-less cell but keeps its boxed style.
$ws2.Range("F8").ClearContents()

# ---------------------------------------------------------------------------
# 3) View state: active tab flips back to "Purchase 22-23", and each sheet's
#    remembered selection moves.
# ---------------------------------------------------------------------------

$ws2.Range("F10").Select()
$ws1.Range("D42").Select()
$ws1.Activate()

$wb.Application.CalculateFull()
